$d = $word.ActiveDocument

$replacements = @(
    @{old="816×4=3264"; new="534×9=4806"},
    @{old="460×8=3680"; new="337×9=3033"},
    @{old="346×6=2076"; new="871×6=5226"},
    @{old="495×8=3960"; new="231×3=693"},
    @{old="702×6=4212"; new="249×3=747"},
    @{old="749×2=1498"; new="112×7=784"},
    @{old="614×3=1842"; new="170×5=850"},
    @{old="265×7=1855"; new="643×4=2572"},
    @{old="838×9=7542"; new="276×4=1104"},
    @{old="916×2=1832"; new="925×4=3700"},
    @{old="845×9=7605"; new="706×4=2824"},
    @{old="180×6=1080"; new="594×3=1782"},
    @{old="682×3=2046"; new="395×7=2765"},
    @{old="583×2=1166"; new="740×5=3700"},
    @{old="758×8=6064"; new="376×8=3008"},
    @{old="587×5=2935"; new="653×8=5224"},
    @{old="723×2=1446"; new="255×2=510"},
    @{old="611×2=1222"; new="977×9=8793"},
    @{old="376×2=752";  new="174×3=522"},
    @{old="774×2=1548"; new="384×8=3072"},
    @{old="372×3=1116"; new="832×4=3328"},
    @{old="382×6=2292"; new="216×8=1728"},
    @{old="906×7=6342"; new="196×3=588"},
    @{old="810×4=3240"; new="107×3=321"},
    @{old="298×6=1788"; new="256×4=1024"}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

$d.Save()
